# "send mail to me every day" — daily report tweak:
#  - hide the working columns F:M (the 1H block) and scroll the view to
#    start at column E so the visible report opens on the summary columns
#  - refresh a handful of stat cells with today's numbers
#  - rename the "万点" header to "sugg_new" and turn the old text
#    "2w"/"1w" suggestion cells into real numeric values
#  - bold the two "I" ratio cells that now stand out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- value refresh -------------------------------------------------------
$ws.Range("Q2").Value = 75
$ws.Range("S2").Value = 73
$ws.Range("V2").Value = 0

$ws.Range("V4").Value = 0.01

$ws.Range("S7").Value = 75
$ws.Range("V7").Value = 0.01

$ws.Range("V13").Value = 0.02

# V22 / V24 used to hold the text suggestions "2w" / "1w" — now plain numbers
$ws.Range("V22").Value = 0.02
$ws.Range("V24").Value = 0

# --- header rename: 万点 -> sugg_new -------------------------------------
$ws.Range("V1").Value = "sugg_new"

# --- bold the two highlighted ratio cells (match the row-22 "call-out" look)
$ws.Range("I22").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- hide columns F:M (the 1H detail block) -------------------------------
$ws.Range("F1:M1").EntireColumn.Hidden = $true

# --- scroll/selection: open the sheet already panned over to column E,
#     with F:M selected (mirrors the hidden block) ------------------------
$ws.Range("F1:M1048576").Select()
$excel.ActiveWindow.ScrollColumn = 5
